# The Employee_Details sheet had a duplicate "county" header column
# (it appeared once at column J and again at column M, with the M
# copy being completely empty below the header). Remove the duplicate
# (empty) "county" column at M so only the original county column
# (J) - which holds real data such as "Dallas" - remains.
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Employee_Details")
$ws.Activate()

$ws.Columns.Item(13).Delete()

$null = $ws.Range("J19").Select()
